$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Change the selection on the "demographics" sheet from M2 to A1:G28
# ------------------------------------------------------------------
$demo = $wb.Worksheets.Item("demographics")
$demo.Activate()
$demo.Range("A1:G28").Select()

# ------------------------------------------------------------------
# 2. Add the new "Sheet1" worksheet at the end of the workbook
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$demoRef = $wb.Worksheets.Item("demographics")
$newSheet.Move($null, $demoRef)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet1"

# ------------------------------------------------------------------
# 3. Column widths for C, D, E
# ------------------------------------------------------------------
$newSheet.Columns.Item(3).ColumnWidth = 19.42578125
$newSheet.Columns.Item(4).ColumnWidth = 12.7109375
$newSheet.Columns.Item(5).ColumnWidth = 14.28515625

# ------------------------------------------------------------------
# 4. Header row
# ------------------------------------------------------------------
$newSheet.Range("A1").Value = "Canton"
$newSheet.Range("B1").Value = "Population"
$newSheet.Range("C1").Value = "SettlementAreaHa"
$newSheet.Range("D1").Value = "SettlementAreaKm2"
$newSheet.Range("E1").Value = "Density"
$newSheet.Range("F1").Value = "O65"
$newSheet.Range("G1").Value = "O65P"
$newSheet.Range("H1").Value = "Beds"
$newSheet.Range("I1").Value = "BedsPerCapita"

# ------------------------------------------------------------------
# 5. Data rows (canton, population, settlement area [ha], O65 share, beds)
# ------------------------------------------------------------------
$data = @(
    @("ZH", 1520968, 37796, 0.17, 4472),
    @("BE", 1034977, 41197, 0.20799999999999999, 3053),
    @("VD", 799145, 29940, 0.16400000000000001, 2268),
    @("AG", 678207, 23854, 0.17699999999999999, 1450),
    @("SG", 507697, 19408, 0.183, 1565),
    @("GE", 495249, 9416, 0.16400000000000001, 1506),
    @("LU", 409557, 14384, 0.17599999999999999, 977),
    @("TI", 353343, 15881, 0.22600000000000001, 1338),
    @("VS", 343955, 18463, 0.19600000000000001, 834),
    @("FR", 318714, 13998, 0.157, 547),
    @("BL", 288132, 9025, 0.219, 582),
    @("TG", 276472, 12170, 0.17599999999999999, 570),
    @("SO", 273194, 10952, 0.19600000000000001, 510),
    @("GR", 198379, 13863, 0.21299999999999999, 546),
    @("BS", 194766, 2628, 0.19900000000000001, 1199),
    @("NE", 177964, 6701, 0.19, 385),
    @("SZ", 159165, 5499, 0.17699999999999999, 274),
    @("ZG", 126837, 3306, 0.17, 206),
    @("SH", 81991, 3403, 0.21199999999999999, 186),
    @("JU", 73419, 5615, 0.20899999999999999, 145),
    @("AR", 55234, 2231, 0.19700000000000001, 208),
    @("NW", 43223, 1481, 0.20300000000000001, 73),
    @("GL", 40403, 1995, 0.20100000000000001, 86),
    @("OW", 37841, 1879, 0.188, 50),
    @("UR", 36433, 2000, 0.20399999999999999, 63),
    @("AI", 16145, 813, 0.191, 18)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Formula = "=C$r*0.01"
    $newSheet.Cells.Item($r, 5).Formula = "=B$r/D$r"
    $newSheet.Cells.Item($r, 6).Value = $row[3]
    $newSheet.Cells.Item($r, 7).Formula = "=ROUND(B$r*F$r,0)"
    $newSheet.Cells.Item($r, 8).Value = $row[4]
    $newSheet.Cells.Item($r, 9).Formula = "=H$r/B$r"
    $r = $r + 1
}

# ------------------------------------------------------------------
# 6. Totals row (28)
# ------------------------------------------------------------------
$newSheet.Range("A28").Value = "CH"
$newSheet.Range("B28").Formula = "=SUM(B2:B27)"
$newSheet.Range("E28").Formula = "=AVERAGE(E2:E27)"
$newSheet.Range("F28").Formula = "=AVERAGE(F2:F27)"
$newSheet.Range("G28").Formula = "=SUM(G2:G27)"
$newSheet.Range("H28").Formula = "=SUM(H2:H27)"
$newSheet.Range("I28").Formula = "=H28/B28"

# ------------------------------------------------------------------
# 7. Sheet view: make this the active/selected tab with its own selection
# ------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("G36").Select()

foreach ($sheet in $wb.Worksheets) {
    Write-Host "Sheet: $($sheet.Name)"
}
